# overwrite excel export --fix
#
# The source export re-ran without the header row ("name" / "username" /
# "plain_password") and with freshly regenerated plain-text passwords for
# every account. Recreate that by deleting the header row (shifting
# everything up) and then writing the new password column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the header row entirely; remaining rows shift up one, and the
# sheet dimension / shared strings are recalculated automatically.
$ws.Rows(1).Delete()

# New plain_password values for each of the 22 remaining rows (column C).
$passwords = @(
    "BnsV29",
    "8sQK33",
    "61eW83",
    "VC0919",
    "0lgY38",
    "GOhc14",
    "qaJc67",
    "s2Z875",
    "CnTP60",
    "NaZ421",
    "LeeY53",
    "QbE547",
    "gLz455",
    "YY4Y11",
    "jCQY95",
    "a9jh80",
    "h3AX24",
    "GGzX42",
    "nuOS65",
    "4gNQ67",
    "BQ7h44",
    "2d0d98"
)

for ($i = 0; $i -lt $passwords.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 3).Value = $passwords[$i]
}
